$d = $word.ActiveDocument

# The document currently ends with a trailing empty paragraph right before
# the sectPr. We replace that paragraph's range with a fresh block of
# paragraphs (new diary entries for April 13/14 2018), finishing with a
# fresh empty paragraph so the document still ends the same way it did
# before the edit.

$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range

$xmlBody = @'
<w:p/><w:p><w:r><w:t>April 13 2018</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Added extra lines to several of the pileups to make them all line up</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Created a named pileup directory</w:t></w:r></w:p><w:p><w:pPr><w:ind w:left="360"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>TODO</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">: make the label names be the file names minus the path and the file </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">extension </w:t></w:r><w:r><w:t>:done</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:r><w:t>April 14, 2018</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Make more space for the names: </w:t></w:r><w:r><w:t>Done increased 80px to 120px</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Build an image for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>consensous</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and variance </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>combind</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/>
'@

$insertRange.InsertXML($xmlBody) | Out-Null
